# Migrated all SAP services to use external ERP service.
#
# The farmer's legacy SAP id ("tusap1") and the two numeric-looking ids that
# used to be stored as real numbers now come back from the external ERP
# system as plain text identifiers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data changes -------------------------------------------------------
# "sopimusmäärä" quantity id on the header row is now ERP-sourced text.
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "100"

# Viljelijän (farmer) SAP id -> new external ERP id.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "111111"

# Marjalajin (berry type) SAP id now comes back from the ERP as text too.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1234"

# --- Formatting that came along with the ERP migration re-save ---------
$rng = $ws.Range("A1:H2")
$rng.Font.Name = "Calibri"
$rng.Font.Size = 12
$rng.Font.Color = 0
$rng.NumberFormat = "General"

$ws.PageSetup.CenterHeader = ""
$ws.PageSetup.CenterFooter = ""

$ws.Range("C3").Select()
